$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Replace the vulnerable DOLLAR() formulas with ABS()
$ws1.Range("H12").Formula = "=ABS(H1)"
$ws2.Range("H1").Formula = "=ABS(Sheet1!F3)"

# Align I1's style with the rest of the row (style de-duplication) by
# copying H1's format (same fill/alignment) onto I1 without touching its formula
$ws2.Range("H1").Copy()
$ws2.Range("I1").PasteSpecial(-4122)

# Restore Sheet1 as the active sheet/tab and update selections
$ws1.Activate()
$ws1.Range("H13").Select() | Out-Null
$ws2.Range("H11").Select() | Out-Null
$ws1.Activate()
